$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'95.900.63"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +0.02%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'3.593.17"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  -2.28%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  +0.00%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'238.29"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  -1.55%  "
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'655.30"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  +1.63%  "
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  +3.83%  "
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  +0.36%  "
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  +0.09%  "
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'1.03"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  +2.24%  "
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'3.594.50"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  -2.25%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'43.10"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  -1.69%  "
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  +1.11%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'6.48"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  +1.12%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'4.261.01"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  -2.55%  "
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.Value = "'95.712.70"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  -0.05%  "
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  -0.60%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'3.577.14"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  -2.81%  "
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'12.80"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  -4.84%  "
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'7.75"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  -2.96%  "
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'18.06"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  -3.67%  "
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'0.496"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  +2.97%  "
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  +0.77%  "
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'512.69"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  -1.33%  "
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'7.09"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +4.26%  "
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  -0.11%  "
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'96.18"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  -1.35%  "
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'12.85"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  +1.35%  "
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'3.783.98"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  -2.20%  "
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "'3.05"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  -3.49%  "
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'0.146"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  +2.87%  "
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  -0.88%  "
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  +0.17%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'0.998"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  -0.43%  "
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  -0.80%  "
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'32.02"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  -3.53%  "
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'1.70"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  +12.83%  "
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  -2.69%  "
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'8.64"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  +8.90%  "
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'598.76"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  +6.46%  "
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  -0.86%  "
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  +0.09%  "
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  +7.02%  "
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  -5.62%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'5.82"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  -0.17%  "
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  +4.30%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'34.50"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  +1.89%  "
$c.Style = "Normal"
$c = $ws.Range("B48")
$c.Value = "'VeChain"
$c.Style = "Normal"
$c = $ws.Range("C48")
$c.Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'0.0420"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -2.86%  "
$c.Style = "Normal"
$c = $ws.Range("B49")
$c.Value = "'WhiteBITCoin"
$c.Style = "Normal"
$c = $ws.Range("C49")
$c.Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'23.47"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  -1.03%  "
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'3.50"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  +0.04%  "
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -0.76%  "
$c.Style = "Normal"
